# Updated cryptos list - apply price/volume changes and two row swaps (32/33, 50/51).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.740.45"
$ws.Range("E2").Value = "  -0.12%  "

$ws.Range("D3").Value = "2.436.95"
$ws.Range("E3").Value = "  -1.23%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "'558.59"

$ws.Range("D6").Value = "'161.55"
$ws.Range("E6").Value = "  -1.07%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("E8").Value = "  -0.13%  "

$ws.Range("E9").Value = "  +7.14%  "

$ws.Range("E10").Value = "  -2.14%  "

$ws.Range("E11").Value = "  -0.63%  "

$ws.Range("D12").Value = "'4.59"
$ws.Range("E12").Value = "  -5.51%  "

$ws.Range("E13").Value = "  +3.43%  "

$ws.Range("D14").Value = "68.630.59"
$ws.Range("E14").Value = "  -0.14%  "

$ws.Range("D15").Value = "2.885.36"
$ws.Range("E15").Value = "  -0.70%  "

$ws.Range("E16").Value = "  -1.54%  "

$ws.Range("D17").Value = "2.437.51"
$ws.Range("E17").Value = "  -0.05%  "

$ws.Range("E18").Value = "  -0.80%  "

$ws.Range("D19").Value = "'339.30"
$ws.Range("E19").Value = "  +0.69%  "

$ws.Range("D20").Value = "'6.94"
$ws.Range("E20").Value = "  -0.03%  "

$ws.Range("E21").Value = "  +1.11%  "

$ws.Range("E22").Value = "  +2.06%  "

$ws.Range("E23").Value = "  -0.04%  "

$ws.Range("D24").Value = "'66.86"
$ws.Range("E24").Value = "  +0.34%  "

$ws.Range("E25").Value = "  +1.42%  "

$ws.Range("D26").Value = "2.564.27"
$ws.Range("E26").Value = "  -1.17%  "

$ws.Range("D27").Value = "'1.00"
$ws.Range("E27").Value = "  +0.43%  "

$ws.Range("E28").Value = "  -0.46%  "

$ws.Range("D29").Value = "0.0₃0820"
$ws.Range("E29").Value = "  -0.34%  "

$ws.Range("E30").Value = "  -1.07%  "

$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").Value = "'1.16"
$ws.Range("E32").Value = "  +1.28%  "

$ws.Range("B33").Value = "Bittensor"
$ws.Range("C33").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D33").Value = "'427.61"
$ws.Range("E33").Value = "  -0.39%  "

$ws.Range("E34").Value = "  -2.06%  "

$ws.Range("E35").Value = "  +0.02%  "

$ws.Range("D36").Value = "'18.98"
$ws.Range("E36").Value = "  -0.15%  "

$ws.Range("E37").Value = "  +0.03%  "

$ws.Range("E38").Value = "  +0.81%  "

$ws.Range("D39").Value = "'0.105"
$ws.Range("E39").Value = "  -3.01%  "

$ws.Range("E40").Value = "  +0.07%  "

$ws.Range("D41").Value = "'1.51"
$ws.Range("E41").Value = "  +2.64%  "

$ws.Range("E42").Value = "  -2.16%  "

$ws.Range("E43").Value = "  -0.68%  "

$ws.Range("D44").Value = "'2.05"
$ws.Range("E44").Value = "  -0.98%  "

$ws.Range("D45").Value = "'130.89"
$ws.Range("E45").Value = "  +0.43%  "

$ws.Range("E46").Value = "  -1.34%  "

$ws.Range("D47").Value = "'0.0719"
$ws.Range("E47").Value = "  +0.00%  "

$ws.Range("E48").Value = "  -0.90%  "

$ws.Range("E49").Value = "  -1.24%  "

$ws.Range("B50").Value = "BitgetToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/q7gMmMdLb+bitgettoken-bgb"
$ws.Range("D50").Value = "'1.15"
$ws.Range("E50").Value = "  +3.26%  "

$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").Value = "'0.0921"
$ws.Range("E51").Value = "  +0.43%  "
